$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wedding seating table headers
$ws.Range("A1").Value = "FIRST NAME"
$ws.Range("B1").Value = "LAST NAME"
$ws.Range("C1").Value = "TABLE NUMBER"

# Row 2
$ws.Range("A2").Value = "JOHN"
$ws.Range("B2").Value = "SILVA"
$ws.Range("C2").Value = 10

# Row 3 (note trailing space kept on first name, as in source data)
$ws.Range("A3").Value = "JOHN "
$ws.Range("B3").Value = "PERERA"
$ws.Range("C3").Value = 11

# Row 4 (table number left blank)
$ws.Range("A4").Value = "JACK"
$ws.Range("B4").Value = "SILVA"

# Size the columns to fit the entered content
$ws.Columns.Item(1).ColumnWidth = 11.71
$ws.Columns.Item(2).ColumnWidth = 11.86
$ws.Columns.Item(3).ColumnWidth = 13.71

# Leave the selection where the author left it
$null = $ws.Range("B10").Select()
